$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new year column (2020), mirroring the
# formatting already used for the preceding year column (P).
$ws.Range("P4:P5").Copy()
$ws.Range("Q4:Q5").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 3.3

# Move the selection, matching the author's saved cursor position.
$ws.Range("R4").Select()
